# Auto-generated edit script applying cryptos list price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '58.899.88'
$ws.Range('E2').Value = '  -0.97%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.501.76'
$ws.Range('E3').Value = '  +2.36%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.22%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '537.41'
$ws.Range('E5').Value = '  +0.49%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '144.07'
$ws.Range('E6').Value = '  -1.96%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.997'

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.570'
$ws.Range('E8').Value = '  +0.42%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.529.53'
$ws.Range('E9').Value = '  +2.81%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0997'
$ws.Range('E10').Value = '  +0.65%  '

# Row 11
$ws.Range('E11').Value = '  +0.18%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '5.53'
$ws.Range('E12').Value = '  +0.93%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.350'
$ws.Range('E13').Value = '  +0.08%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.953.44'
$ws.Range('E14').Value = '  +2.45%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '23.52'
$ws.Range('E15').Value = '  -2.42%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '58.823.84'
$ws.Range('E16').Value = '  -0.99%  '

# Row 17
$ws.Range('E17').Value = '  +0.81%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.528.32'
$ws.Range('E18').Value = '  +1.47%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.20'
$ws.Range('E19').Value = '  +0.20%  '

# Row 20
$ws.Range('E20').Value = '  -2.57%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '322.86'
$ws.Range('E21').Value = '  -0.62%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').Value = '  +2.80%  '

# Row 23
$ws.Range('E23').Value = '  +1.05%  '

# Row 24
$ws.Range('E24').Value = '  +2.14%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.437'
$ws.Range('E25').Value = '  -6.31%  '

# Row 26
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.163'
$ws.Range('E26').Value = '  +0.63%  '

# Row 27
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.629.13'
$ws.Range('E27').Value = '  +2.68%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  +2.26%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.76'
$ws.Range('E29').Value = '  +0.46%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.76'
$ws.Range('E30').Value = '  -1.65%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0₃0771'
$ws.Range('E31').Value = '  +0.87%  '

# Row 32
$ws.Range('E32').Value = '  -0.94%  '

# Row 33
$ws.Range('E33').Value = '  -7.38%  '

# Row 34
$ws.Range('E34').Value = '  -0.16%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '157.94'
$ws.Range('E35').Value = '  +1.16%  '

# Row 36
$ws.Range('E36').Value = '  +6.73%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '18.62'
$ws.Range('E37').Value = '  +1.49%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.36'
$ws.Range('E38').Value = '  -3.86%  '

# Row 39
$ws.Range('E39').Value = '  -7.68%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.57'
$ws.Range('E40').Value = '  -3.56%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '36.48'
$ws.Range('E41').Value = '  -1.04%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '299.55'
$ws.Range('E42').Value = '  -4.76%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.65'
$ws.Range('E43').Value = '  -1.97%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.815'
$ws.Range('E44').Value = '  -4.22%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.995'
$ws.Range('E45').Value = '  -0.21%  '

# Row 46
$ws.Range('E46').Value = '  +3.86%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.77'
$ws.Range('E47').Value = '  +0.43%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '125.46'
$ws.Range('E48').Value = '  +4.71%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0929'
$ws.Range('E49').Value = '  -0.96%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '18.65'
$ws.Range('E50').Value = '  +0.29%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0514'
$ws.Range('E51').Value = '  -2.22%  '
